# Adds a "Categoria" column (Perros / Gatos / Chivos) between the
# "Numero" and "Hora de salida" columns, bumps Juana's "Numero" from
# 12 to 14, and resets the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Hora de salida" used to live in column C; it slides over to column D.
# Recreate it there (value + formatting) before column C is repurposed.
$ws.Cells.Item(1, 4).Font.Name = "Arial"
$ws.Cells.Item(1, 4).Font.Size = 10
$ws.Cells.Item(1, 4).Value = $ws.Cells.Item(1, 3).Value2

# "Hora de salida" moved from ~08:29 to 00:30 for every entry (stored as
# the Excel day-fraction serial for that time).
$horaSalida = @{
    2  = 0.0208333333333333
    3  = 0.0208333333333333
    4  = 0.0208333333333334
    5  = 0.0208333333333334
    6  = 0.0208333333333334
    7  = 0.0208333333333335
    8  = 0.0208333333333335
    9  = 0.0208333333333335
    10 = 0.0208333333333336
    11 = 0.0208333333333336
    12 = 0.0208333333333336
    13 = 0.0208333333333336
    14 = 0.0208333333333337
    15 = 0.0208333333333337
    16 = 0.0208333333333337
    17 = 0.0208333333333338
    18 = 0.0208333333333338
    19 = 0.0208333333333338
}

foreach ($row in $horaSalida.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.NumberFormat = "HH:MM:SS"
    $cell.Value = $horaSalida[$row]
}

# Widen the "Hora de salida" column slightly; everything else keeps the
# sheet's default width.
$ws.Columns.Item(4).ColumnWidth = 12.63

# Column C becomes "Categoria".
$ws.Range("C1").Value = "Categoria"

# Category assignments, by row group.
$categories = @{
    2  = "Perros"
    3  = "Perros"
    4  = "Perros"
    5  = "Perros"
    6  = "Perros"
    7  = "Gatos"
    8  = "Gatos"
    9  = "Gatos"
    10 = "Gatos"
    11 = "Gatos"
    12 = "Gatos"
    13 = "Gatos"
    14 = "Chivos"
    15 = "Chivos"
    16 = "Chivos"
    17 = "Chivos"
    18 = "Chivos"
    19 = "Chivos"
}

foreach ($row in $categories.Keys) {
    $cell = $ws.Cells.Item($row, 3)
    # Column C used to hold the time-formatted "Hora de salida" values;
    # drop that formatting now that it holds plain category text (back to
    # the sheet's plain default look, same as column B's data cells).
    $cell.NumberFormat = "General"
    $cell.Font.Name = "Noto Sans CJK SC"
    $cell.Font.Size = 10
    $cell.Value = $categories[$row]
}

# Juana's "Numero" changed from 12 to 14.
$ws.Range("B2").Value = 14

# Restore the previously-selected cell to the new location recorded in
# the workbook (F12 instead of F18).
$ws.Range("F12").Select()
